$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# parameters sheet: elasticity value A2 changes from 0.2 to 0.02, and the
# selection moves off this sheet (it stops being the "tabSelected" sheet).
# ---------------------------------------------------------------------------
$wsParams = $wb.Worksheets.Item("parameters")
$wsParams.Range("A2").Value = 0.02
$wsParams.Range("A3").Select()

# ---------------------------------------------------------------------------
# tech sheet: a few technology numbers change, and the selection moves.
# ---------------------------------------------------------------------------
$wsTech = $wb.Worksheets.Item("tech")
$wsTech.Range("G3").Value = 0
$wsTech.Range("D4").Value = 8
$wsTech.Range("G4").Value = 10
$wsTech.Range("I8").Select()

# ---------------------------------------------------------------------------
# day_weights sheet: weights are rebalanced, and this becomes the active
# (selected) tab/sheet, with B2:B4 selected (active cell B2).
# ---------------------------------------------------------------------------
$wsDayWeights = $wb.Worksheets.Item("day_weights")
$wsDayWeights.Range("B2").Value = 199
$wsDayWeights.Range("B3").Value = 106
$wsDayWeights.Range("B4").Value = 61

# ---------------------------------------------------------------------------
# cap_factors sheet: the capacity-factor profiles for summer, fall_spring and
# winter (rows 2-4) are recomputed; update every changed cell with the new
# literal values from the recalculated model. Selection also changes.
# ---------------------------------------------------------------------------
$wsCapFactors = $wb.Worksheets.Item("cap_factors")

# Row 2 - summer
$wsCapFactors.Range("G2").Value = 4.0201005025125598E-3
$wsCapFactors.Range("H2").Value = 4.1170854271356763E-2
$wsCapFactors.Range("I2").Value = 0.15498492462311561
$wsCapFactors.Range("J2").Value = 0.34120100502512568
$wsCapFactors.Range("K2").Value = 0.51137185929648243
$wsCapFactors.Range("L2").Value = 0.63440201005025132
$wsCapFactors.Range("M2").Value = 0.70589447236180902
$wsCapFactors.Range("N2").Value = 0.7173165829145729
$wsCapFactors.Range("O2").Value = 0.67734170854271358
$wsCapFactors.Range("P2").Value = 0.57809547738693468
$wsCapFactors.Range("Q2").Value = 0.43412562814070349
$wsCapFactors.Range("R2").Value = 0.25252763819095481
$wsCapFactors.Range("S2").Value = 9.1608040201005048E-2
$wsCapFactors.Range("T2").Value = 1.7597989949748739E-2
$wsCapFactors.Range("U2").Value = 3.2160804020100472E-4

# Row 3 - fall_spring
$wsCapFactors.Range("G3").Value = 2.1509433962264152E-3
$wsCapFactors.Range("H3").Value = 2.0773584905660381E-2
$wsCapFactors.Range("I3").Value = 7.4056603773584917E-2
$wsCapFactors.Range("J3").Value = 0.20321698113207551
$wsCapFactors.Range("K3").Value = 0.34733018867924531
$wsCapFactors.Range("L3").Value = 0.46702830188679251
$wsCapFactors.Range("M3").Value = 0.54499999999999993
$wsCapFactors.Range("N3").Value = 0.55942452830188683
$wsCapFactors.Range("O3").Value = 0.51176415094339622
$wsCapFactors.Range("P3").Value = 0.41073584905660382
$wsCapFactors.Range("Q3").Value = 0.27456603773584909
$wsCapFactors.Range("R3").Value = 0.13354716981132081
$wsCapFactors.Range("S3").Value = 3.995283018867922E-2
$wsCapFactors.Range("T3").Value = 8.6698113207547139E-3
$wsCapFactors.Range("U3").Value = 1.4150943396226421E-4

# Row 4 - winter
$wsCapFactors.Range("G4").Value = 1.8032786885245899E-3
$wsCapFactors.Range("H4").Value = 1.6032786885245912E-2
$wsCapFactors.Range("I4").Value = 5.4229508196721322E-2
$wsCapFactors.Range("J4").Value = 0.1274918032786885
$wsCapFactors.Range("K4").Value = 0.20950819672131141
$wsCapFactors.Range("L4").Value = 0.26809836065573772
$wsCapFactors.Range("M4").Value = 0.28267213114754092
$wsCapFactors.Range("N4").Value = 0.27595081967213131
$wsCapFactors.Range("O4").Value = 0.25488524590163941
$wsCapFactors.Range("P4").Value = 0.21809836065573759
$wsCapFactors.Range("Q4").Value = 0.151
$wsCapFactors.Range("R4").Value = 8.1081967213114725E-2
$wsCapFactors.Range("S4").Value = 3.1426229508196707E-2
$wsCapFactors.Range("T4").Value = 7.7049180327868824E-3
$wsCapFactors.Range("U4").Value = 9.8360655737705021E-5

$wsCapFactors.Range("B3:Y3").Select()

# ---------------------------------------------------------------------------
# Finally, activate day_weights so it becomes the workbook's active sheet
# (tabSelected) with its B2:B4 selection intact, matching the saved view
# state in the target workbook.
# ---------------------------------------------------------------------------
$wsDayWeights.Activate()
$wsDayWeights.Range("B2:B4").Select()
